# Add MetadataSheet to all Templates
$wb = $excel.ActiveWorkbook
$excel.UserName = "Oliver Maus"

# --- 1. Rename existing sheet, add the new metadata sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "4COM04_GenomeAssembly"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "SwateTemplateMetadata"

# --- 2. Cell content (Col A = field labels, Col B/C/D = values) ---
$data = @(
  @(1,1,"Id"), @(1,2,"11b23480-80a8-4d95-a2bb-b3e9c7d53a23"),
  @(2,1,"Name"), @(2,2,"Genome assembly"),
  @(3,1,"Version"), @(3,2,"1.1.3"),
  @(4,1,"Description"), @(4,2,"Template to describe computational of a genome assembly"),
  @(5,1,"Docslink"),
  @(6,1,"Organisation"),
  @(7,1,"Table"), @(7,2,"annotationTableSmoothBird74"),
  @(8,1,"#ER list"),
  @(9,1,"ER"), @(9,2,"SRA"), @(9,3,"GENBANK"),
  @(10,1,"ER Term Accession Number"),
  @(11,1,"ER Term Source REF"),
  @(12,1,"#TAGS list"),
  @(13,1,"Tags"), @(13,2,"Genomics"), @(13,3,"DNASeq"), @(13,4,"Assay"),
  @(14,1,"Tags Term Accession Number"),
  @(15,1,"Tags Term Source REF"),
  @(16,1,"#AUTHORS list"),
  @(17,1,"Authors Last Name"), @(17,2,"Kranz"), @(17,3,"Brilhaus"), @(17,4,"Maus"),
  @(18,1,"Authors First Name"), @(18,2,"Angela"), @(18,3,"Dominik"), @(18,4,"Oliver"),
  @(19,1,"Authors Mid Initials"),
  @(20,1,"Authors Email"),
  @(21,1,"Authors Phone"),
  @(22,1,"Authors Fax"),
  @(23,1,"Authors Address"),
  @(24,1,"Authors Affiliation"),
  @(25,1,"#AUTHORS ROLES list"),
  @(26,1,"Authors Roles"),
  @(27,1,"Authors Roles Term Accession Number"),
  @(28,1,"Authors Roles Term Source REF")
)

foreach ($entry in $data) {
  $r = $entry[0]
  $c = $entry[1]
  $v = $entry[2]
  $ws2.Cells.Item($r, $c).Value = $v
}

Write-Host "content written"
